$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "# cities" header in H1, matching the style used by the rest
# of the header row (bold, centered, bordered).
$ws.Range("H1").Value = "# cities"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the "# cities" values for each state, rows 2-33.
$values = @(2022,5545,2543,2762,21157,12186,4034,1226,634,5890,8809,6769,4690,10348,4894,8644,1578,2850,4822,10723,6568,2192,2207,6554,5495,7300,2472,6566,1175,19845,2434,4498)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
